# Split the single "Hypothyroidism" sheet into two sheets:
#   - "HypothyroidismAdd"       : recipes to add (kept short / trimmed)
#   - "HypothyroidismEliminate" : full original recipe list (untouched copy)

$wb = $excel.ActiveWorkbook

# The existing sheet becomes the "Add" sheet.
$wsAdd = $wb.Worksheets.Item(1)
$wsAdd.Name = "HypothyroidismAdd"

# Create the "Eliminate" sheet right after it and copy the full original data
# (all 52 recipe rows + header) into it before trimming the "Add" sheet.
$wsEliminate = $wb.Worksheets.Add($null, $wsAdd)
$wsEliminate.Name = "HypothyroidismEliminate"

$fullRange = $wsAdd.Range("A1:K52")
$fullRange.Copy($wsEliminate.Range("A1"))

# The header row only spans columns A-J; the copy above created a spurious
# empty K1 cell because the source range is rectangular - remove it again.
$wsEliminate.Range("K1").ClearContents()

# On the "Add" sheet, the Masala Chawli recipe (row 2) only needs the base
# ingredient to be added to the diet - shorten its ingredients list.
$wsAdd.Range("E2").Value = "1/2 cup chawli (cow pea / lobhia) , soaked overnight and drained    salt to taste    "

# The "Add" sheet only keeps the first two recipes (rows 2-3); drop the rest.
$wsAdd.Range("A4:K52").ClearContents()
